$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value2 = 567
$ws.Range("C4").Value2 = '上海·元宵AuPoRo音乐动漫FES（取消）'
$ws.Range("D4").Value2 = '友谊路街道友谊路318号 灏唯滨江智创园'
$ws.Range("E4").Value2 = '2024.02.24 08:00-02.24 21:00'
$ws.Range("F4").Value2 = 26
$ws.Range("G4").Value2 = '不可售'
$ws.Range("H4").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81481'
$ws.Range("I4").Value2 = '//i1.hdslb.com/bfs/openplatform/202401/w4kr4a0X1706497345456.jpeg'
$ws.Range("C5").Value2 = '上海·原X铁X崩only（取消）'
$ws.Range("D5").Value2 = '澳门路168号 月星国际家居'
$ws.Range("E5").Value2 = '2024.02.24 10:30-02.24 16:30'
$ws.Range("F5").Value2 = 173
$ws.Range("H5").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81446'
$ws.Range("I5").Value2 = '//i2.hdslb.com/bfs/openplatform/202401/IIePRulM1706248855263.jpeg'
$ws.Range("C6").Value2 = '上海·原神×崩坏×星铁only旅行盛宴2.0'
$ws.Range("D6").Value2 = '西藏南路1号 上海大世界'
$ws.Range("E6").Value2 = '2024.02.24 10:00-02.25 17:00'
$ws.Range("F6").Value2 = 3147
$ws.Range("G6").Value2 = 65
$ws.Range("H6").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81276'
$ws.Range("I6").Value2 = '//i2.hdslb.com/bfs/openplatform/202401/82hU3z8m1706155835021.png'
$ws.Range("C7").Value2 = '上海·第三届燃梦BACG国潮嘉年华-原X铁X崩同好交流'
$ws.Range("D7").Value2 = '盈浦街道淀山浦社区淀山湖大道851号青浦万达茂F3 万达汽车乐园(青浦万达茂店)'
$ws.Range("E7").Value2 = '2024.02.24 11:00-02.24 16:30'
$ws.Range("F7").Value2 = 2733
$ws.Range("G7").Value2 = 65.8
$ws.Range("H7").Value2 = 'https://show.bilibili.com/platform/detail.html?id=77754'
$ws.Range("I7").Value2 = '//i0.hdslb.com/bfs/openplatform/202402/JYUdM9Q91707963393893.jpeg'
$ws.Range("C8").Value2 = '上海·趣元界&斗罗大陆上元佳节次元派对'
$ws.Range("D8").Value2 = '长宁路1191号长宁来福士B1 长宁来福士'
$ws.Range("E8").Value2 = '2024.02.24 11:30-02.25 17:30'
$ws.Range("F8").Value2 = 530
$ws.Range("G8").Value2 = 98
$ws.Range("H8").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81415'
$ws.Range("I8").Value2 = '//i0.hdslb.com/bfs/openplatform/202401/yis4JHfE1706169986733.jpeg'
$ws.Range("C9").Value2 = '上海·魔都元宵节漫展-COS为王'
$ws.Range("D9").Value2 = '澳门路168号月星家居(江宁路地铁站1号口步行420米) 月星广场'
$ws.Range("E9").Value2 = '2024.02.24 10:00-02.25 16:00'
$ws.Range("F9").Value2 = 42
$ws.Range("G9").Value2 = 49
$ws.Range("H9").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81238'
$ws.Range("I9").Value2 = '//i1.hdslb.com/bfs/openplatform/202401/KxQZPADR1705913896609.jpeg'
$ws.Range("C10").Value2 = '上海·魔都多厨狂喜漫展-CH01'
$ws.Range("F10").Value2 = 18
$ws.Range("H10").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81423'
$ws.Range("I10").Value2 = '//i1.hdslb.com/bfs/openplatform/202401/axpOY3zo1706173660010.jpeg'
$ws.Range("B11").Value2 = '2024.03.02'
$ws.Range("C11").Value2 = '上海·原神X星穹铁道ONLY'
$ws.Range("D11").Value2 = '逸仙路301号靠纪念路路口 上海宝丰联大酒店'
$ws.Range("E11").Value2 = '2024.03.02 10:00-03.02 17:00'
$ws.Range("F11").Value2 = 343
$ws.Range("G11").Value2 = 60
$ws.Range("H11").Value2 = 'https://show.bilibili.com/platform/detail.html?id=80299'
$ws.Range("I11").Value2 = '//i2.hdslb.com/bfs/openplatform/202312/V0xu26Cl1703753850690.jpeg'
$ws.Range("B12").Value2 = '2024.03.03'
$ws.Range("C12").Value2 = '上海·怀旧番ONLY'
$ws.Range("D12").Value2 = '逸仙路270号  上海宝丰联大酒店'
$ws.Range("E12").Value2 = '2024.03.03 10:00-03.03 17:00'
$ws.Range("F12").Value2 = 282
$ws.Range("H12").Value2 = 'https://show.bilibili.com/platform/detail.html?id=80575'
$ws.Range("I12").Value2 = '//i1.hdslb.com/bfs/openplatform/202401/y4uWdyPT1704700763902.jpeg'
$ws.Range("B13").Value2 = '2024.03.08'
$ws.Range("C13").Value2 = '上海·第八届ACBC动漫盛典-国潮汉服游园会'
$ws.Range("D13").Value2 = '浦锦南路1586弄2号 奇迹花园'
$ws.Range("E13").Value2 = '2024.03.08 10:00-03.10 17:00'
$ws.Range("F13").Value2 = 26
$ws.Range("H13").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81456'
$ws.Range("I13").Value2 = '//i1.hdslb.com/bfs/openplatform/202401/qZtpawf51706254849667.jpeg'
$ws.Range("B14").Value2 = '2024.03.09'
$ws.Range("C14").Value2 = '上海·S·CGE动漫游戏嘉年华'
$ws.Range("D14").Value2 = '军工路1076号 纪希片场(秀场)'
$ws.Range("E14").Value2 = '2024.03.09 10:00-03.10 17:00'
$ws.Range("F14").Value2 = 5626
$ws.Range("G14").Value2 = 70
$ws.Range("H14").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81173'
$ws.Range("I14").Value2 = '//i0.hdslb.com/bfs/openplatform/202401/TYA5FLkE1705891815532.jpeg'
$ws.Range("C15").Value2 = '上海·第五十三届燃梦星辰动漫嘉年华-随机宅舞'
$ws.Range("D15").Value2 = '周家嘴路3608号 宝龙旭辉广场'
$ws.Range("E15").Value2 = '2024.03.09 10:20-03.10 16:30'
$ws.Range("F15").Value2 = 610
$ws.Range("G15").Value2 = 58
$ws.Range("H15").Value2 = 'https://show.bilibili.com/platform/detail.html?id=80571'
$ws.Range("I15").Value2 = '//i0.hdslb.com/bfs/openplatform/202401/SHH70VXN1704700240858.jpeg'
$ws.Range("C16").Value2 = '上海·青山刚昌ONLY【名侦探柯南&魔术快斗】'
$ws.Range("D16").Value2 = '漕宝路1688号 诺宝中心酒店'
$ws.Range("E16").Value2 = '2024.03.09 10:00-03.09 16:30'
$ws.Range("F16").Value2 = 1017
$ws.Range("G16").Value2 = 73
$ws.Range("H16").Value2 = 'https://show.bilibili.com/platform/detail.html?id=76410'
$ws.Range("I16").Value2 = '//i2.hdslb.com/bfs/openplatform/202309/fVXMrcHy1693971682397.jpeg'
$ws.Range("B17").Value2 = '2024.03.16'
$ws.Range("C17").Value2 = '上海·Look Look动漫嘉年华'
$ws.Range("D17").Value2 = '龙吴路4800号2号门 有只怪兽片场'
$ws.Range("E17").Value2 = '2024.03.16 10:00-03.17 17:30'
$ws.Range("F17").Value2 = 53
$ws.Range("G17").Value2 = 29.9
$ws.Range("H17").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81804'
$ws.Range("I17").Value2 = '//i2.hdslb.com/bfs/openplatform/202402/WFRql6sg1707274094000.jpeg'
$ws.Range("C18").Value2 = '上海·SISP动漫游戏嘉年华'
$ws.Range("D18").Value2 = '年家浜路518号 周浦万达广场'
$ws.Range("E18").Value2 = '2024.03.16 13:00-03.17 19:00'
$ws.Range("F18").Value2 = 165
$ws.Range("G18").Value2 = 48
$ws.Range("H18").Value2 = 'https://show.bilibili.com/platform/detail.html?id=80339'
$ws.Range("I18").Value2 = '//i0.hdslb.com/bfs/openplatform/202312/a8iuOufB1703832570508.jpeg'
$ws.Range("F19").Value2 = 80
$ws.Range("F20").Value2 = 447
$ws.Range("F21").Value2 = 1215
$ws.Range("F23").Value2 = 110
$ws.Range("F24").Value2 = 328
$ws.Range("F25").Value2 = 42
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value2 = 240
$ws.Range("F9").Value2 = 51
$ws.Range("F13").Value2 = 629
$ws.Range("F18").Value2 = 52
$ws.Range("F23").Value2 = 336
$ws.Range("F25").Value2 = 4030
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value2 = 2551
$ws.Range("F6").Value2 = 1109
$ws.Range("F9").Value2 = 1429
$ws.Range("F10").Value2 = 398
$ws.Range("F11").Value2 = 112
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value2 = 2551
$ws.Range("F6").Value2 = 1109
$ws.Range("F7").Value2 = 1429
$ws.Range("F8").Value2 = 398
$ws.Range("F9").Value2 = 112
$ws.Range("F12").Value2 = 567
$ws.Range("C13").Value2 = '上海·《哈利的魔法世界》动漫视听音乐会'
$ws.Range("D13").Value2 = '都市路4889号（莘庄地铁站南广场） 上海保利城市剧院'
$ws.Range("E13").Value2 = '2024.02.24 14:30-02.24 16:00'
$ws.Range("F13").Value2 = 22
$ws.Range("G13").Value2 = 158
$ws.Range("H13").Value2 = 'https://show.bilibili.com/platform/detail.html?id=80639'
$ws.Range("I13").Value2 = '//i2.hdslb.com/bfs/openplatform/202401/4PieCC9N1706261750579.jpeg'
$ws.Range("C14").Value2 = '上海·原神×崩坏×星铁only旅行盛宴2.0'
$ws.Range("D14").Value2 = '西藏南路1号 上海大世界'
$ws.Range("E14").Value2 = '2024.02.24 10:00-02.25 17:00'
$ws.Range("F14").Value2 = 3147
$ws.Range("G14").Value2 = 65
$ws.Range("H14").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81276'
$ws.Range("I14").Value2 = '//i2.hdslb.com/bfs/openplatform/202401/82hU3z8m1706155835021.png'
$ws.Range("C15").Value2 = '上海·第三届燃梦BACG国潮嘉年华-原X铁X崩同好交流'
$ws.Range("D15").Value2 = '盈浦街道淀山浦社区淀山湖大道851号青浦万达茂F3 万达汽车乐园(青浦万达茂店)'
$ws.Range("E15").Value2 = '2024.02.24 11:00-02.24 16:30'
$ws.Range("F15").Value2 = 2733
$ws.Range("G15").Value2 = 65.8
$ws.Range("H15").Value2 = 'https://show.bilibili.com/platform/detail.html?id=77754'
$ws.Range("I15").Value2 = '//i0.hdslb.com/bfs/openplatform/202402/JYUdM9Q91707963393893.jpeg'
$ws.Range("C16").Value2 = '上海·趣元界&斗罗大陆上元佳节次元派对'
$ws.Range("D16").Value2 = '长宁路1191号长宁来福士B1 长宁来福士'
$ws.Range("E16").Value2 = '2024.02.24 11:30-02.25 17:30'
$ws.Range("F16").Value2 = 530
$ws.Range("G16").Value2 = 98
$ws.Range("H16").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81415'
$ws.Range("I16").Value2 = '//i0.hdslb.com/bfs/openplatform/202401/yis4JHfE1706169986733.jpeg'
$ws.Range("C17").Value2 = '上海·魔都元宵节漫展-COS为王'
$ws.Range("D17").Value2 = '澳门路168号月星家居(江宁路地铁站1号口步行420米) 月星广场'
$ws.Range("E17").Value2 = '2024.02.24 10:00-02.25 16:00'
$ws.Range("F17").Value2 = 42
$ws.Range("G17").Value2 = 49
$ws.Range("H17").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81238'
$ws.Range("I17").Value2 = '//i1.hdslb.com/bfs/openplatform/202401/KxQZPADR1705913896609.jpeg'
$ws.Range("C18").Value2 = '上海·魔都多厨狂喜漫展-CH01'
$ws.Range("F18").Value2 = 18
$ws.Range("H18").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81423'
$ws.Range("I18").Value2 = '//i1.hdslb.com/bfs/openplatform/202401/axpOY3zo1706173660010.jpeg'
$ws.Range("B19").Value2 = '2024.02.25'
$ws.Range("C19").Value2 = '上海·青山吉能见面会'
$ws.Range("D19").Value2 = '虹许路731号4号楼 THE BOXX•城市乐园'
$ws.Range("E19").Value2 = '2024.02.25 14:30-02.25 19:30'
$ws.Range("F19").Value2 = 240
$ws.Range("G19").Value2 = 380
$ws.Range("H19").Value2 = 'https://show.bilibili.com/platform/detail.html?id=80142'
$ws.Range("I19").Value2 = '//i0.hdslb.com/bfs/openplatform/202312/1npuHFBM1703231674558.jpeg'
$ws.Range("B20").Value2 = '2024.03.02'
$ws.Range("C20").Value2 = '上海·2024藤田玲上海粉丝见面会'
$ws.Range("D20").Value2 = '宜昌路179号 万代南梦宫上海文化中心'
$ws.Range("E20").Value2 = '2024.03.02 12:30-03.02 19:40'
$ws.Range("F20").Value2 = 22
$ws.Range("G20").Value2 = 580
$ws.Range("H20").Value2 = 'https://show.bilibili.com/platform/detail.html?id=80993'
$ws.Range("I20").Value2 = '//i0.hdslb.com/bfs/openplatform/202401/Vm6ntgVd1705548188785.png'
$ws.Range("C21").Value2 = '上海·原神X星穹铁道ONLY'
$ws.Range("D21").Value2 = '逸仙路301号靠纪念路路口 上海宝丰联大酒店'
$ws.Range("E21").Value2 = '2024.03.02 10:00-03.02 17:00'
$ws.Range("F21").Value2 = 343
$ws.Range("G21").Value2 = 60
$ws.Range("H21").Value2 = 'https://show.bilibili.com/platform/detail.html?id=80299'
$ws.Range("I21").Value2 = '//i2.hdslb.com/bfs/openplatform/202312/V0xu26Cl1703753850690.jpeg'
$ws.Range("C22").Value2 = '上海·小山百代2024上海粉丝见面会'
$ws.Range("D22").Value2 = '宜昌路179号 万代南梦宫上海文化中心'
$ws.Range("E22").Value2 = '2024.03.02 13:00-03.02 20:00'
$ws.Range("F22").Value2 = 334
$ws.Range("G22").Value2 = 380
$ws.Range("H22").Value2 = 'https://show.bilibili.com/platform/detail.html?id=80924'
$ws.Range("I22").Value2 = '//i1.hdslb.com/bfs/openplatform/202401/FpA9OkKy1705467080070.jpeg'
$ws.Range("B23").Value2 = '2024.03.03'
$ws.Range("C23").Value2 = '上海·“前方核能！！！”和音社交响乐团·经典动漫音乐会'
$ws.Range("D23").Value2 = '丁香路425号 上海东方艺术中心'
$ws.Range("E23").Value2 = '2024.03.03 19:30-03.03 21:00'
$ws.Range("F23").Value2 = 51
$ws.Range("G23").Value2 = 162
$ws.Range("H23").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81844'
$ws.Range("I23").Value2 = '//i2.hdslb.com/bfs/openplatform/202402/eWndSczF1707386523895.jpeg'
$ws.Range("B24").Value2 = '2024.03.08'
$ws.Range("C24").Value2 = '上海·第八届ACBC动漫盛典-国潮汉服游园会'
$ws.Range("D24").Value2 = '浦锦南路1586弄2号 奇迹花园'
$ws.Range("E24").Value2 = '2024.03.08 10:00-03.10 17:00'
$ws.Range("F24").Value2 = 26
$ws.Range("G24").Value2 = 60
$ws.Range("H24").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81456'
$ws.Range("I24").Value2 = '//i1.hdslb.com/bfs/openplatform/202401/qZtpawf51706254849667.jpeg'
$ws.Range("B25").Value2 = '2024.03.09'
$ws.Range("C25").Value2 = '上海·S·CGE动漫游戏嘉年华'
$ws.Range("D25").Value2 = '军工路1076号 纪希片场(秀场)'
$ws.Range("E25").Value2 = '2024.03.09 10:00-03.10 17:00'
$ws.Range("F25").Value2 = 5626
$ws.Range("G25").Value2 = 70
$ws.Range("H25").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81173'
$ws.Range("I25").Value2 = '//i0.hdslb.com/bfs/openplatform/202401/TYA5FLkE1705891815532.jpeg'
$ws.Range("C26").Value2 = '上海·《挪威的森林》—摇滚情歌之夜演唱会'
$ws.Range("D26").Value2 = '南京西路1376号 上海商城剧院'
$ws.Range("E26").Value2 = '2024.03.09 19:30-03.09 21:00'
$ws.Range("F26").Value2 = 1
$ws.Range("G26").Value2 = 72
$ws.Range("H26").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81241'
$ws.Range("I26").Value2 = '//i2.hdslb.com/bfs/openplatform/202401/1FJ0Fj5m1705915336335.jpeg'
$ws.Range("C27").Value2 = '上海·爱乐之城音乐会'
$ws.Range("E27").Value2 = '2024.03.09 14:00-03.09 15:30'
$ws.Range("F27").Value2 = 13
$ws.Range("G27").Value2 = 108
$ws.Range("H27").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81289'
$ws.Range("I27").Value2 = '//i2.hdslb.com/bfs/openplatform/202401/ZZXtDrwZ1705996679699.jpeg'
$ws.Range("C28").Value2 = '上海·第五十三届燃梦星辰动漫嘉年华-随机宅舞'
$ws.Range("D28").Value2 = '周家嘴路3608号 宝龙旭辉广场'
$ws.Range("E28").Value2 = '2024.03.09 10:20-03.10 16:30'
$ws.Range("F28").Value2 = 610
$ws.Range("G28").Value2 = 58
$ws.Range("H28").Value2 = 'https://show.bilibili.com/platform/detail.html?id=80571'
$ws.Range("I28").Value2 = '//i0.hdslb.com/bfs/openplatform/202401/SHH70VXN1704700240858.jpeg'
$ws.Range("C29").Value2 = '上海·青山刚昌ONLY【名侦探柯南&魔术快斗】'
$ws.Range("D29").Value2 = '漕宝路1688号 诺宝中心酒店'
$ws.Range("E29").Value2 = '2024.03.09 10:00-03.09 16:30'
$ws.Range("F29").Value2 = 1017
$ws.Range("G29").Value2 = 73
$ws.Range("H29").Value2 = 'https://show.bilibili.com/platform/detail.html?id=76410'
$ws.Range("I29").Value2 = '//i2.hdslb.com/bfs/openplatform/202309/fVXMrcHy1693971682397.jpeg'
$ws.Range("B30").Value2 = '2024.03.10'
$ws.Range("C30").Value2 = '上海·三森铃子10周年纪念2024演唱会'
$ws.Range("D30").Value2 = '宜昌路179号 万代南梦宫上海文化中心'
$ws.Range("E30").Value2 = '2024.03.10 18:00-03.10 19:30'
$ws.Range("F30").Value2 = 629
$ws.Range("G30").Value2 = 399
$ws.Range("H30").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81433'
$ws.Range("I30").Value2 = '//i0.hdslb.com/bfs/openplatform/202401/L8rmm2h81706236781799.jpeg'
$ws.Range("B31").Value2 = '2024.03.16'
$ws.Range("C31").Value2 = '上海·Look Look动漫嘉年华'
$ws.Range("D31").Value2 = '龙吴路4800号2号门 有只怪兽片场'
$ws.Range("E31").Value2 = '2024.03.16 10:00-03.17 17:30'
$ws.Range("F31").Value2 = 53
$ws.Range("G31").Value2 = 29.9
$ws.Range("H31").Value2 = 'https://show.bilibili.com/platform/detail.html?id=81804'
$ws.Range("I31").Value2 = '//i2.hdslb.com/bfs/openplatform/202402/WFRql6sg1707274094000.jpeg'
$ws.Range("C32").Value2 = '上海·SISP动漫游戏嘉年华'
$ws.Range("D32").Value2 = '年家浜路518号 周浦万达广场'
$ws.Range("E32").Value2 = '2024.03.16 13:00-03.17 19:00'
$ws.Range("F32").Value2 = 165
$ws.Range("G32").Value2 = 48
$ws.Range("H32").Value2 = 'https://show.bilibili.com/platform/detail.html?id=80339'
$ws.Range("I32").Value2 = '//i0.hdslb.com/bfs/openplatform/202312/a8iuOufB1703832570508.jpeg'
$ws.Range("F33").Value2 = 80
$ws.Range("F36").Value2 = 52
$ws.Range("F39").Value2 = 336
$ws.Range("F40").Value2 = 1215
$ws.Range("F47").Value2 = 328
$ws.Range("F48").Value2 = 42
